$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "41.337.48"
$ws.Cells.Item(2, 5).Value = "  -0.78%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.437.11"
$ws.Cells.Item(3, 5).Value = "  -1.51%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.21%  "

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "316.63"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.47%  "

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "89.54"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -3.21%  "

# Row 7
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.542"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -2.11%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.16%  "

# Row 9
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.496"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -3.74%  "

# Row 10
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "32.07"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -2.54%  "

# Row 11
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.0827"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -8.12%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -2.55%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "2.813.89"
$ws.Cells.Item(13, 5).Value = "  -1.45%  "

# Row 14
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "6.69"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -3.08%  "

# Row 15
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "15.29"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -2.40%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "2.432.08"
$ws.Cells.Item(16, 5).Value = "  -2.65%  "

# Row 17
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "0.770"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -2.60%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "41.235.42"
$ws.Cells.Item(18, 5).Value = "  -0.92%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "0.0₃0919"
$ws.Cells.Item(19, 5).Value = "  -4.55%  "

# Row 20
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "6.23"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -3.65%  "

# Row 21
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "71.67"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +0.39%  "

# Row 22
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "11.04"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -3.45%  "

# Row 23
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "234.76"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -2.72%  "

# Row 24
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "2.69"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -1.97%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +0.13%  "

# Row 26
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "1.88"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -2.17%  "

# Row 27
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "23.96"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -3.65%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -3.30%  "

# Row 29
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "9.52"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -3.41%  "

# Row 30
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "34.62"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -4.49%  "

# Row 31
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "157.37"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +0.49%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +0.11%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "Filecoin"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "5.25"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -4.82%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -1.45%  "

# Row 35
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "0.0742"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -3.19%  "

# Row 36
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "2.90"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -0.46%  "

# Row 37
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "16.53"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -5.14%  "

# Row 38
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.114"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -1.03%  "

# Row 39
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "1.77"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -3.46%  "

# Row 40
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "0.0996"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -3.33%  "

# Row 41
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "3.88"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -2.48%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  -7.16%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "1.982.97"
$ws.Cells.Item(43, 5).Value = "  +0.15%  "

# Row 44
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.0274"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -3.70%  "

# Row 45
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "18.04"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -5.73%  "

# Row 46
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "2.85"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -4.59%  "

# Row 47
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "9.44"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +3.10%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "2.675.67"
$ws.Cells.Item(48, 5).Value = "  -1.22%  "

# Row 49
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "95.11"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -2.33%  "

# Row 50
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "73.00"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -1.07%  "

# Row 51
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "51.72"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -1.56%  "
